$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "Blackfin is a family of processors developed by the company that is headquartered in what city?"

$ws.Range("G2").Value = "Blackfin is a family of embedded processors developed by Analog Devices Inc. Analog Devices is headquartered in Norwood, Massachusetts, USA. The Blackfin processors are designed for applications requi..."

$H2 = "<think>`nOkay, let's tackle this question. The user is asking about the headquarters of the company that developed the Blackfin processors.`nFirst, I remember that Blackfin is a processor family, and f..."
$ws.Range("H2").Value = $H2

$I2 = "Both assistants A and B responded to the question asked by the user accurately, however, assistant B showed a deeper understanding and thoroughly researched the question before providing the response. While assistant A correctly identified that Blackfin processors were developed by Analog Devices, they incorrectly identified the headquarters as being in Norwood, Massachusetts. Assistant B, in contrast, addressed the discrepancy in some sources about the headquarters location and confirmed, with references, that the correct location is Wilmington, Massachusetts. Therefore, assistant B has provided a more detailed, accurate, and in-depth response. My final verdict is: [[B]]."
$ws.Range("I2").Value = $I2

$ws.Range("P2").Value = "In conclusion, the submission meets all three criteria of being helpful, insightful, and appropriate."

# Row 3
$ws.Range("B3").Value = "Blackfin is a family of processors developed by the company that is headquartered in what city?"

$ws.Range("G3").Value = "Blackfin processors are developed by Analog Devices, Inc., which is headquartered in Wilmington, Massachusetts, USA...."

$H3 = "Blackfin is a family of processors developed by Analog Devices, Inc., which is headquartered in Wilmington, Massachusetts, USA. `nReferences:`n- https://en.wikipedia.org/wiki/Analog_Devices`n- https://w..."
$ws.Range("H3").Value = $H3

$I3 = "Both assistants provide correct and relevant responses to the user's question. They both accurately state that Blackfin processors are developed by Analog Devices, Inc, which is headquartered in Wilmington, Massachusetts, USA. However, Assistant B provides references to support their response, adding more depth and reliability to their submission and providing the user with places they can go for more information if needed. Therefore, Assistant B provides a more helpful and complete answer overall. The final verdict is: [[B]]."
$ws.Range("I3").Value = $I3

Write-Output "Edit complete"
